$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "season record" columns (Wins / Losses / Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the exact same bold / centered / bordered header style as the rest of row 1
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in Cleveland's 2017 season record (102 Wins, 60 Losses, 0 Ties)
# for every player data row (2 through 43)
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 102  # AD: Wins
    $ws.Cells.Item($row, 31).Value = 60   # AE: Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF: Ties
}
